# edit.ps1 - applies the changes described by the XML diff to before.docx
#
# Summary of changes:
#  1) In the "Környezet" section, the run "...írtam Jav" + a _GoBack bookmark +
#     the run "a 11-es verzión...JavaFX-et." are merged back into a single
#     continuous sentence "...írtam Java 11-es verzión...JavaFX-et." (and the
#     _GoBack bookmark is removed from this location).
#  2) In the installation instructions paragraph, the sentence describing the
#     XAMPP/PhpMyAdmin database import is reworded: the clause describing
#     creating the 'timetable' database is dropped, and "beimportálni" (to
#     import) is moved earlier in the sentence, right after "kell".
#  3) After the run containing ".Main)" a new run containing just a line
#     break (<w:br/>) is appended within the same paragraph.
#  4) The _GoBack bookmark re-appears, this time at the start of the
#     paragraph that only contains two line breaks (near the end of the
#     document, just before the last image anchor paragraph).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: merge the two runs split by the _GoBack bookmark back into one
# continuous run, removing the bookmark from this spot.
# ---------------------------------------------------------------------
$old1 = "Az adatbázis az XAMPP nevő programból futtattam, az adatbáziskezelő " + `
        "természetesen MySQL. Az alkalmazást Java nyelven írtam Java 11-es " + `
        "verzión. Adatbázishoz való csatlakozáshoz JDBC-t használtam, a gui " + `
        "összerakásához pedig JavaFX-et."
$new1 = $old1
$rng = $d.Content
$found1 = $rng.Find.Execute("Az adatbázis az XAMPP nevő programból futtattam, az adatbáziskezelő természetesen MySQL. Az alkalmazást Java nyelven írtam Jav" + "a 11-es verzión. Adatbázishoz való csatlakozáshoz JDBC-t használtam, a gui összerakásához pedig JavaFX-et.", `
                             $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Edit1 (merge sentence / drop bookmark):" $found1

# ---------------------------------------------------------------------
# Edit 2: reword the XAMPP/PhpMyAdmin installation sentence.
# ---------------------------------------------------------------------
$quote = [char]0x2019
$old2 = "Az adatbázis előkészítéséhez XAMPP-ban kell MySQL szervert futtatni, majd PhpMyAdmin-ban kell létrehozni a " + `
        $quote + "timetable" + $quote + " nevű adatbázist, amibe a projektben található " + `
        $quote + "timetable.sql" + $quote + " fájlt kell beimportálni."
$new2 = "Az adatbázis előkészítéséhez XAMPP-ban kell MySQL szervert futtatni, majd PhpMyAdmin-ban kell beimportálni a projektben található " + `
        $quote + "timetable.sql" + $quote + " fájlt."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Host "Edit2 (reword XAMPP sentence):" $found2

# ---------------------------------------------------------------------
# Edit 3: append a new run with just a line break after ".Main)".
# ---------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(".Main)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Edit3 find .Main):" $found3
if ($found3) {
    $rng3.Collapse(0)
    $rng3.InsertAfter([char]0x000B)
}

# ---------------------------------------------------------------------
# Edit 4: re-insert the _GoBack bookmark at the start of the paragraph
# that contains only two line breaks (near the end of the document).
# ---------------------------------------------------------------------
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Az Órarendek tabon van lehetőségünk", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Edit4 find anchor paragraph:" $found4
if ($found4) {
    $rng4.Expand(4)                                  # wdParagraph -> whole paragraph
    $nextParaRange = $d.Range($rng4.End, $rng4.End)
    $nextParaRange.Expand(4)                          # the following (two-line-break) paragraph
    $bmPos = $nextParaRange.Start
    $d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
    Write-Host "Bookmark added at" $bmPos
}

Write-Host "Done."
